$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.788.16'
$ws.Range("E2").Value = '  -2.61%  '

$ws.Range("D3").Value = '1.744.93'
$ws.Range("E3").Value = '  -5.04%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  +0.05%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '238.62'
$ws.Range("E5").Value = '  -8.89%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9998'
$ws.Range("E6").Value = '  -0.02%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5046'
$ws.Range("E7").Value = '  -6.30%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '42.01'
$ws.Range("E8").Value = '  -6.31%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2729'
$ws.Range("E9").Value = '  -9.51%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.06154'
$ws.Range("E10").Value = '  -10.77%  '

$ws.Range("D11").Value = '1.746.11'
$ws.Range("E11").Value = '  -5.05%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.06927'
$ws.Range("E12").Value = '  -2.98%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '15.51'
$ws.Range("E13").Value = '  -12.15%  '

$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.530'
$ws.Range("E14").Value = '  -9.30%  '

$ws.Range("B15").Value = 'Polygon'
$ws.Range("C15").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6022'
$ws.Range("E15").Value = '  -18.45%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '76.98'
$ws.Range("E16").Value = '  -13.66%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.000'
$ws.Range("E17").Value = '  -0.02%  '

$ws.Range("E18").Value = '  -0.03%  '

$ws.Range("D19").Value = '25.792.19'
$ws.Range("E19").Value = '  -2.70%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000006889'
$ws.Range("E20").Value = '  -12.80%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.62'
$ws.Range("E21").Value = '  -16.10%  '

$ws.Range("D22").Value = '1.968.37'
$ws.Range("E22").Value = '  -5.28%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.061'
$ws.Range("E23").Value = '  -11.58%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.241'
$ws.Range("E24").Value = '  -12.42%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '8.171'
$ws.Range("E25").Value = '  -11.26%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '137.93'
$ws.Range("E26").Value = '  -3.42%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.466'
$ws.Range("E27").Value = '  -14.69%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.818'
$ws.Range("E28").Value = '  -16.70%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.00'
$ws.Range("E29").Value = '  -11.84%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '103.86'
$ws.Range("E30").Value = '  -6.59%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08136'
$ws.Range("E31").Value = '  -8.11%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.708'
$ws.Range("E32").Value = '  -12.74%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.484'
$ws.Range("E33").Value = '  -14.02%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.04546'
$ws.Range("E34").Value = '  -6.12%  '

$ws.Range("E35").Value = '  -0.02%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.615'
$ws.Range("E36").Value = '  -10.69%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.9861'
$ws.Range("E37").Value = '  -12.89%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.6106'
$ws.Range("E38").Value = '  -16.38%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.682'
$ws.Range("E39").Value = '  -13.40%  '

$ws.Range("E40").Value = '  -9.51%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.937'
$ws.Range("E41").Value = '  -14.40%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.9994'
$ws.Range("E42").Value = '  -0.02%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '101.92'

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.3848'
$ws.Range("E44").Value = '  -18.36%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.7381'
$ws.Range("E45").Value = '  -18.42%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '4.974'
$ws.Range("E46").Value = '  -15.71%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.05378'
$ws.Range("E47").Value = '  -6.73%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.1114'
$ws.Range("E48").Value = '  -11.08%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '5.973'
$ws.Range("E49").Value = '  -19.33%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '30.20'

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '52.57'
$ws.Range("E51").Value = '  -12.54%  '
